$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the Confessarius character from "赦罪師" to "聴罪師" throughout the
# Japanese dialogue column (B), including the inline mention within the
# "？？？" speaker's line in row 42.
$rows = @(42, 44, 45, 46, 47, 48, 50, 52, 53, 54, 56, 57, 59, 60)
foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value()
    if ($val -ne $null) {
        $newVal = $val.Replace("赦罪師", "聴罪師")
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}

# Update the English localization quote style in C60 from curly quotes to
# straight quotes: "She" is ready. -> 'She' is ready.
$c60 = $ws.Cells.Item(60, 3)
$c60Val = $c60.Value()
if ($c60Val -ne $null) {
    $openQuote = [string]([char]0x201C)
    $closeQuote = [string]([char]0x201D)
    $newC60 = $c60Val.Replace($openQuote + "She" + $closeQuote, "'She'")
    if ($newC60 -ne $c60Val) {
        $c60.Value = $newC60
    }
}
